# Update Name of Algo
# Applies updated KNN-imputed values to the affected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B7"   = 5.355
    "A9"   = -21.743
    "B12"  = 5.57
    "D13"  = -7.726999999999999
    "C15"  = -13.247
    "D16"  = -8.529
    "A18"  = -22.051
    "A20"  = -20.793
    "D20"  = -7.558
    "D24"  = -7.547
    "B26"  = 5.492
    "A27"  = -21.446
    "B27"  = 6.603
    "E27"  = 16.741
    "B29"  = 5.360999999999999
    "E29"  = 17.123
    "B37"  = 8.670999999999999
    "B38"  = 5.038
    "C38"  = -12.566
    "D39"  = -7.507
    "C44"  = -12.248
    "D48"  = -7.366000000000001
    "B51"  = 5.790999999999999
    "C51"  = -11.915
    "D52"  = -7.811
    "B55"  = 5.755
    "D56"  = -7.825
    "C57"  = -13.222
    "E57"  = 16.557
    "C63"  = -12.076
    "A69"  = -21.375
    "B69"  = 6.165000000000001
    "B70"  = 5.606
    "C70"  = -11.206
    "A76"  = -20.718
    "A82"  = -22.025
    "B83"  = 6.010000000000001
    "D84"  = -8.181000000000001
    "E85"  = 16.651
    "C99"  = -12.396
    "D100" = -8.310999999999998
    "D101" = -7.831
    "B102" = 7.398999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
